# Insert a new data row at row 254 (pushing the existing rows 254-314 down to
# 255-315) and populate the new row with a new "Zanahoria" price record for
# "Vega Modelo de Temuco" (weekly fruit/vegetable price update).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 254..314 down to 255..315, keeping formatting (e.g. the date
# style on column D) of the row being pushed down.
$ws.Rows.Item(254).Insert()

# Populate the newly-inserted row 254 with the new record.
$ws.Cells.Item(254, 1).Value = 10
$ws.Cells.Item(254, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(254, 3).Value = "La Araucanía"
$ws.Cells.Item(254, 4).Value = 44722
$ws.Cells.Item(254, 5).Value = 9
$ws.Cells.Item(254, 6).Value = 100114013
$ws.Cells.Item(254, 7).Value = "Zanahoria"
$ws.Cells.Item(254, 8).Value = "Sin especificar"
$ws.Cells.Item(254, 9).Value = "Primera"
$ws.Cells.Item(254, 10).Value = 80
$ws.Cells.Item(254, 11).Value = 6000
$ws.Cells.Item(254, 12).Value = 6000
$ws.Cells.Item(254, 13).Value = 6000
$ws.Cells.Item(254, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(254, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(254, 16).Value = 240
$ws.Cells.Item(254, 17).Value = 25
$ws.Cells.Item(254, 18).Value = "Hortaliza"
